# "Add files via upload / updates" --------------------------------------
#
# The diff fills in a block of previously-blank quarterly growth cells on
# Sheet1 (AE5:AH5 and AE6:AH6 - the "qtrs-travellers"-derived % change
# formulas that mirror the pattern already used by the surrounding
# columns), and flips which sheet/cell is the active selection: Sheet1
# becomes the tab-selected sheet (selection AE6) instead of
# qtrs-travellers (which reverts to its own stored selection, I12).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$qtrs   = $wb.Worksheets.Item("qtrs-travellers")

# --- Fill in the previously-empty percentage-change formulas --------------
# Row 5: (row3 - row2) / row2   (2019 vs 2018 growth)
$sheet1.Range("AE5").Formula = "=(AE3-AE2)/AE2"
$sheet1.Range("AF5").Formula = "=(AF3-AF2)/AF2"
$sheet1.Range("AG5").Formula = "=(AG3-AG2)/AG2"
$sheet1.Range("AH5").Formula = "=(AH3-AH2)/AH2"

# Row 6: (row4 - row3) / row3   (2020 vs 2019 growth)
$sheet1.Range("AE6").Formula = "=(AE4-AE3)/AE3"
$sheet1.Range("AF6").Formula = "=(AF4-AF3)/AF3"
$sheet1.Range("AG6").Formula = "=(AG4-AG3)/AG3"
$sheet1.Range("AH6").Formula = "=(AH4-AH3)/AH3"

# --- Restore qtrs-travellers' own selection before switching away --------
$qtrs.Activate()
$qtrs.Range("I12").Select()

# --- Make Sheet1 the active/tab-selected sheet with the new selection ----
$sheet1.Activate()
$sheet1.Range("AE6").Select()
